$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Series")

# --- Update row 132 (was the placeholder "Labor Force Status Flows" series) to the
#     LNS17400000 / "Labor Force Flows Employed to Unemployed" series ---
$ws2.Range("C132").Value = "LNS17400000"
$ws2.Range("D132").Value = "Labor Force Flows Employed to Unemployed"
$ws2.Range("F132").Value = "m_change"

# --- Update row 133 (was the placeholder QCEW series) to the QCEWEMP / "Total Employment" series ---
$ws2.Range("C133").Value = "QCEWEMP"
$ws2.Range("D133").Value = "Total Employment"
$ws2.Range("F133").Value = "m_growth"

# --- New row 135: QCEW Average Weekly Wage series ---
$ws2.Range("A135").Value = "Quarterly Census of Employment and Wages"
$ws2.Range("B135").Value = 362
$ws2.Range("C135").Value = "QCEWWAGE"
$ws2.Range("D135").Value = "Average Weekly Wage"
$ws2.Range("E135").Value = "US Bureau of Labor Statistics"
$ws2.Range("F135").Value = "q_growth"
$ws2.Range("G135").Value = "Labor Market"

# --- Sheet1: register the new release (id 180 - Unemployment Insurance Weekly Claims Report) ---
$ws1.Range("A47").Value = 180
$ws1.Range("B47").Value = "Unemployment Insurance Weekly Claims Report"

# --- New row 136: Initial Unemployment Claims series ---
$ws2.Range("A136").Value = "Unemployment Insurance Weekly Claims Report"
$ws2.Range("B136").Value = 180
$ws2.Range("D136").Value = "Initial Unemployment Claims"
$ws2.Range("C136").Value = "IC4WSA"
$ws2.Range("E136").Value = "US Employment and Training Administration"
$ws2.Range("F136").Value = "w_growth"
$ws2.Range("G136").Value = "Labor Market"

# --- Restore the selection/active-cell state shown in the final workbook ---
$ws1.Range("B47").Select()
$ws2.Range("D136").Select()
